$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)
$shp = $s.Shapes.Item(7)
$tr = $shp.TextFrame.TextRange

# Locate the run of text we need to rewrite: "класс Path." -> "класc System.IO.Path."
$full = $tr.Text
$oldChunk = "класс Path."
$idx0 = $full.IndexOf($oldChunk)
$start = $idx0 + 1

# Replace "класс Path." with the new plain (non-monospace / non-teal) text.
# This both strips the special "Path" formatting (teal Courier New) and
# sets the final wording in one shot, merging the affected runs into a
# single plain run that inherits the plain formatting of the text
# immediately before it ("и папок можно применять ").
$chunk = $tr.Characters($start, $oldChunk.Length)
$chunk.Text = "класc System.IO.Path."

# Now re-split that merged run into the individual pieces so the run
# boundaries match the authored edit: "клас" | "c " | "System.IO.Path" | "."
$rc = $tr.Characters($start + 4, 2)
$rc.Text = "c "

$rDot = $tr.Characters($start + 4 + 2 + 14, 1)
$rDot.Text = "."
